$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("AssignmentTable")
$ws.Range("E3").Value = "BDH Emergency (PGY1/2)"
$ws.Range("F3").Value = "CMN Emergency (PGY1/2)"
$ws.Range("G3").Value = "JHH Emergency (PGY1/2)"
$ws.Range("H3").Value = "MH Emergency (PGY1/2)"
$ws = $wb.Worksheets.Item("EvaluationTable")
$ws.Range("B4").Value = "'100.00"
$ws.Range("B6").Value = "'3.00"
$ws.Range("B7").Value = "'5.00"
$ws.Range("B8").Value = "'1.00"
$ws.Range("B9").Value = "'1.58"
$ws = $wb.Worksheets.Item("DoctorPreferenceReport")
$ws.Range("B8").Value = "BDH Emergency (PGY1/2)"
$ws.Range("C8").Value = 2
$ws.Range("B9").Value = "CMN Emergency (PGY1/2)"
$ws.Range("C9").Value = 3
$ws.Range("B10").Value = "JHH Emergency (PGY1/2)"
$ws.Range("C10").Value = 4
$ws.Range("B11").Value = "MH Emergency (PGY1/2)"
$ws.Range("C11").Value = 1
$ws = $wb.Worksheets.Item("HETIComplianceTable")
$ws.Range("B2").Value = "'3"
$ws.Range("B3").Value = "'3"
$ws.Range("D3").Value = "'FALSE"
$ws.Range("F3").Value = "'FALSE"
$ws = $wb.Worksheets.Item("UnassignedTerms")
$ws.Range("B5").Value = "BDH Sub-Acute General Medicine (Team 5) (PGY1/2)"
$ws.Range("C5").Value = "BDH Emergency (PGY1/2)"
$ws.Range("B6").Value = "BDH Surgery (PGY1/2)"
$ws.Range("C6").Value = "BDH Surgery (PGY1/2)"
$ws.Range("B7").Value = "CMN Emergency (PGY1/2)"
$ws.Range("B10").Value = "CMN General Medicine (PGY1/2)"
$ws.Range("B14").Value = "CMN Medicine - Gastroenterology with Drug & Alcohol (PGY1/2)"
$ws.Range("B15").Value = "CMN Medicine - Haematology (01) (PGY1/2)"
$ws.Range("B16").Value = "CMN Medicine - Medical Oncology (01) (PGY1/2)"
$ws.Range("B17").Value = "CMN Palliative Care (PGY1/2)"
$ws.Range("B18").Value = "CMN Surgery (PGY1/2)"
$ws.Range("B20").Value = "HNE Mental Health & Substance Use Service - Mater Hospital (PGY1/2)"
$ws.Range("B21").Value = "HNE Mental Health Lake Macquarie Ward - Mater hospital (PGY1/2)"
$ws.Range("B22").Value = "HNE Mental Health Newcastle Inpatient Unit - Mater Hospital (PGY1/2)"
$ws.Range("B23").Value = "HNE Mental Health Older People's Service - Mater Hospital (PGY1/2)"
$ws.Range("B24").Value = "JHH Acute General Surgery Unit (PGY1/2)"
$ws.Range("B26").Value = "JHH Emergency (PGY1/2)"
$ws.Range("D33").Value = "JHH General Medicine (PGY1/2)"
$ws.Range("B34").Value = "JHH General Medicine (PGY1/2)"
$ws.Range("D37").Value = "JHH General Surgery / Colorectal (PGY1/2)"
$ws.Range("E37").Value = "JHH General Surgery / Hepatopancreatobiliary and Upper GI Surgery (PGY1/2)"
$ws.Range("B38").Value = "JHH General Surgery / Colorectal (PGY1/2)"
$ws.Range("D38").Value = "JHH General Surgery / Trauma (PGY1/2)"
$ws.Range("E38").Value = "JHH General Surgery / Trauma (PGY1/2)"
$ws.Range("B39").Value = "JHH General Surgery / Trauma (PGY1/2)"
$ws.Range("D39").Value = "JHH Medicine - Cardiology (PGY1/2)"
$ws.Range("B40").Value = "JHH General Surgery / Upper GI (PGY1/2)"
$ws.Range("D40").Value = "JHH Medicine - Dermatology (PGY1/2)"
$ws.Range("E40").Value = "JHH Medicine - Cardiology (PGY1/2)"
$ws.Range("B41").Value = "JHH Medicine - Cardiology (PGY1/2)"
$ws.Range("D41").Value = "JHH Medicine - Gastroenterology (PGY1/2)"
$ws.Range("B42").Value = "JHH Medicine - Dermatology (PGY1/2)"
$ws.Range("D42").Value = "JHH Medicine - MACU (PGY1/2)"
$ws.Range("E42").Value = "JHH Medicine - Gastroenterology (PGY1/2)"
$ws.Range("B43").Value = "JHH Medicine - Gastroenterology (PGY1/2)"
$ws.Range("D43").Value = "JHH Medicine - Nephrology & Dialysis (PGY1/2)"
$ws.Range("E43").Value = "JHH Medicine - Infectious Disease (PGY1/2)"
$ws.Range("B44").Value = "JHH Medicine - MACU (PGY1/2)"
$ws.Range("E44").Value = "JHH Medicine - Neurology (PGY1/2)"
$ws.Range("B45").Value = "JHH Medicine - Nephrology & Dialysis (PGY1/2)"
$ws.Range("B46").Value = "JHH Medicine - Neurology (PGY1/2)"
$ws.Range("E46").Value = "JHH Medicine - Rehabilitation (PGY1/2)"
$ws.Range("B47").Value = "JHH Medicine - Rehabilitation (PGY1/2)"
$ws.Range("E47").Value = "JHH Medicine - Respiratory (PGY1/2)"
$ws.Range("E48").Value = "JHH Obstetrics & Gynaecology (PGY1/2)"
$ws.Range("B49").Value = "JHH Medicine - Respiratory (PGY1/2)"
$ws.Range("E49").Value = "JHH Surgery - Cardiothoracic Surgery (PGY1/2)"
$ws.Range("B51").Value = "JHH Obstetrics & Gynaecology (PGY1/2)"
$ws.Range("E51").Value = "JHH Surgery - Neurosurgery (PGY1/2)"
$ws.Range("E56").Value = "JHH Surgery - Orthopaedics (PGY1/2)"
$ws.Range("E57").Value = "JHH Surgery - Orthopaedics (Team 6) (PGY1/2)"
$ws.Range("E59").Value = "JHH Surgery - Urology (PGY1/2)"
$ws.Range("E61").Value = "JHH Surgery - Vascular (PGY1/2)"
$ws.Range("E62").Value = "MH  Medicine - Cardiology (PGY1/2)"
